$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.991.33"
$ws.Range("E2").Value = "  +2.85%  "

$ws.Range("D3").Value = "1.864.60"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.01"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6393"
$ws.Range("E6").Value = "  +3.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2999"
$ws.Range("E8").Value = "  +3.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07495"
$ws.Range("E9").Value = "  +2.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  +6.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07688"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").Value = "1.870.62"
$ws.Range("E12").Value = "  +3.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.068"
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6908"
$ws.Range("E14").Value = "  +4.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.31"
$ws.Range("E15").Value = "  +2.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009469"
$ws.Range("E16").Value = "  +6.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.108"
$ws.Range("E17").Value = "  +4.59%  "

$ws.Range("D18").Value = "29.961.03"
$ws.Range("E18").Value = "  +2.82%  "

$ws.Range("D19").Value = "2.122.42"
$ws.Range("E19").Value = "  +2.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.41"
$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.71"
$ws.Range("E21").Value = "  +2.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.472"
$ws.Range("E23").Value = "  +4.43%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.37"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1425"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.594"
$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.02"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06145"
$ws.Range("E29").Value = "  +10.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.508"
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.274"
$ws.Range("E31").Value = "  +5.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.140"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.153"
$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.879"
$ws.Range("E34").Value = "  +3.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.169"
$ws.Range("E35").Value = "  +3.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7317"
$ws.Range("E36").Value = "  -0.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.603"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.868"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01805"
$ws.Range("E39").Value = "  +2.74%  "

$ws.Range("D40").Value = "1.225.57"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9312"
$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.279"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.025.47"
$ws.Range("E44").Value = "  +3.08%  "

$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.55"
$ws.Range("E46").Value = "  +2.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5091"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.337"
$ws.Range("E48").Value = "  +3.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("E49").Value = "  -8.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4105"
$ws.Range("E50").Value = "  +2.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1142"
$ws.Range("E51").Value = "  +3.09%  "
